$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New menu row added by the user via the xlsx import (id / title / description)
$ws.Range("A19").Value = "05360824-639e-471c-a44f-127064d32a98"
$ws.Range("B19").Value = "новое меню"
$ws.Range("C19").Value = "меню нво"

# Leave the selection on the newly entered cell, matching the saved view state
$ws.Range("C19").Select()

# Page setup as configured before saving
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
